$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, $targetText, $startIndex) {
    for ($i = $startIndex; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Insert a new bullet "Find a Web service for historical exchange
#    rates" right before the "Medium" heading paragraph.
#    We clone the paragraph-level and run-level formatting from an
#    existing bullet that already has the desired look (ListParagraph
#    style, ilvl=1/numId=3, justified, HTMLCode character style with
#    the minor-theme fonts at size 11pt/22 half-points), then simply
#    overwrite its cloned text.
# ------------------------------------------------------------------

$donorIndex  = Get-ParagraphIndexByText $d "Check indexes for all the tables and add if necessary" 1
$formatDonor = $d.Paragraphs.Item($donorIndex)
$donorRange  = $d.Range($formatDonor.Range.Start, $formatDonor.Range.End)

$mediumIndex    = Get-ParagraphIndexByText $d "Medium" 1
$mediumHeading  = $d.Paragraphs.Item($mediumIndex)

$insertPoint = $d.Range($mediumHeading.Range.Start, $mediumHeading.Range.Start)
$insertPoint.FormattedText = $donorRange.FormattedText

# The clone was inserted immediately before "Medium", i.e. it now
# occupies the paragraph slot that used to belong to "Medium". Find it
# by searching *starting at* that same slot so the original donor
# paragraph (earlier in the document) is not matched instead.
$newParaIndex = Get-ParagraphIndexByText $d "Check indexes for all the tables and add if necessary" $mediumIndex
$newPara = $d.Paragraphs.Item($newParaIndex)
$newTextRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newTextRange.Text = "Find a Web service for historical exchange rates"

# ------------------------------------------------------------------
# 2) Remove the old "Get Exchange rate from web service" bullet - its
#    job is now done by the paragraph inserted above.
# ------------------------------------------------------------------

$oldIndex = Get-ParagraphIndexByText $d "Get Exchange rate from web service" 1
$oldPara = $d.Paragraphs.Item($oldIndex)
$oldPara.Range.Delete()
